$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance (week 9, column K) for the students present in this lab session
$ws.Range("K13").Value = $true
$ws.Range("K18").Value = $true
$ws.Range("J42").Value = $true
$ws.Range("K43").Value = $true

# Update the active selection to reflect where the user ended up after editing
$ws.Range("K44").Select()
